# Update NMA and MA coefficient tables to use M instead of mu, for
# consistency with the PDF documentation.
#
# Every "mu_<n>" text value in column H (the mu_name column) across all
# four sheets becomes "M_<n>" - same numeric suffix, just a renamed
# prefix. We do this with a plain find/replace on each cell's value
# rather than hard-coding positions, so it's robust either way.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    # Row 1 is the "mu_name" header - leave it alone; only the data
    # values below it ("mu_1", "mu_2", ...) get renamed.
    foreach ($r in 2..$rows) {
        $cell = $ws.Cells.Item($r, 8)  # column H = mu_name
        $val = $cell.Value2
        if ($val -ne $null -and $val.ToString().StartsWith("mu_")) {
            $cell.Value2 = "M_" + $val.ToString().Substring(3)
        }
    }
}

# The workbook was left with the "weibull" sheet active (rather than
# "fracpoly2"), and on each sheet the selection now sits one row below
# the last data row in column H.
$wsWeibull = $wb.Worksheets.Item("weibull")
$wsWeibull.Activate() | Out-Null
$wsWeibull.Range("H8").Select() | Out-Null

$wb.Worksheets.Item("gompertz").Range("H8").Select() | Out-Null
$wb.Worksheets.Item("fracpoly1").Range("H10").Select() | Out-Null
$wb.Worksheets.Item("fracpoly2").Range("H10").Select() | Out-Null

$wsWeibull.Activate() | Out-Null
